$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 132
$ws.Range("J132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H132").Value = 1977.25
$ws.Range("L132").Value = 0
$ws.Range("I132").Value = 1977.25
$ws.Range("M132").Value = -3401.75
$ws.Range("K132").Value = 5931.75
# row 133
$ws.Range("H133").Value = 12553390
$ws.Range("L133").Value = 12553390
$ws.Range("J133").Value = 12553390
$ws.Range("N133").Value = -12563510
# row 135
$ws.Range("M135").Value = -45001587
$ws.Range("H135").Value = 3334228
$ws.Range("I135").Value = 5000458
$ws.Range("K135").Value = 45004122
# row 138
$ws.Range("I138").Value = 2076.8965
$ws.Range("K138").Value = 6230.689499999999
$ws.Range("H138").Value = 1616801
$ws.Range("N138").Value = -9117683
$ws.Range("J138").Value = 3035801
$ws.Range("L138").Value = 9107403
$ws.Range("M138").Value = -1090.689499999999

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("K32").Value = 4171108.2
$ws.Range("H32").Value = 3927817.5
$ws.Range("I32").Value = 4171108.2
$ws.Range("M32").Value = -4170821.2
# row 57
$ws.Range("K57").Value = 4999.143
$ws.Range("I57").Value = 4999.143
$ws.Range("H57").Value = 4999.143
$ws.Range("M57").Value = -4515.143
# row 74
$ws.Range("I74").Value = 47327.273
$ws.Range("M74").Value = -46453.273
$ws.Range("L74").Value = 4442.3335
$ws.Range("H74").Value = 32191.412
$ws.Range("J74").Value = 4442.3335
$ws.Range("K74").Value = 47327.273
$ws.Range("N74").Value = -6190.3335
# row 77
$ws.Range("J77").Value = 4442.3335
$ws.Range("M77").Value = -232268.365
$ws.Range("H77").Value = 32191.412
$ws.Range("K77").Value = 236636.365
$ws.Range("N77").Value = -30947.6675
$ws.Range("I77").Value = 47327.273
$ws.Range("L77").Value = 22211.6675
# row 97
$ws.Range("K97").Value = 3469.5454
$ws.Range("L97").Value = 13892541
$ws.Range("I97").Value = 3469.5454
$ws.Range("N97").Value = -13893533
$ws.Range("H97").Value = 4905495
$ws.Range("M97").Value = -2973.5454
$ws.Range("J97").Value = 13892541
# row 102
$ws.Range("J102").Value = 4511
$ws.Range("I102").Value = 2701.5334
$ws.Range("M102").Value = -1079.5334
$ws.Range("K102").Value = 2701.5334
$ws.Range("L102").Value = 4511
$ws.Range("N102").Value = -7755
$ws.Range("H102").Value = 3380.0833
# row 122
$ws.Range("M122").Value = -2240.7502
$ws.Range("K122").Value = 4690.7502
$ws.Range("H122").Value = 3955.24
$ws.Range("I122").Value = 1563.5834
# row 126
$ws.Range("H126").Value = 4399.143
$ws.Range("M126").Value = -10727.429
$ws.Range("K126").Value = 13197.429
$ws.Range("I126").Value = 4399.143
# row 132
$ws.Range("H132").Value = 6016.067
$ws.Range("I132").Value = 4471.0464
$ws.Range("M132").Value = -10883.1392
$ws.Range("K132").Value = 13413.1392

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 26
$ws.Range("I26").Value = 29996
$ws.Range("H26").Value = 37978
$ws.Range("M26").Value = -29704
$ws.Range("K26").Value = 29996
# row 64
$ws.Range("L64").Value = 1723.45
$ws.Range("J64").Value = 1723.45
$ws.Range("N64").Value = -2173.45
$ws.Range("H64").Value = 10418072
# row 67
$ws.Range("H67").Value = 10418072
$ws.Range("J67").Value = 1723.45
$ws.Range("N67").Value = -3283.45
$ws.Range("L67").Value = 1723.45
# row 94
$ws.Range("J94").Value = 6931.6665
$ws.Range("H94").Value = 2261.2917
$ws.Range("N94").Value = -7833.6665
$ws.Range("L94").Value = 6931.6665
# row 105
$ws.Range("H105").Value = 6189.6577
$ws.Range("J105").Value = 5134.067
$ws.Range("L105").Value = 5134.067
$ws.Range("N105").Value = -8628.066999999999
# row 107
$ws.Range("M107").Value = -75005720
$ws.Range("I107").Value = 75007640
$ws.Range("K107").Value = 75007640
$ws.Range("H107").Value = 66183570
# row 113
$ws.Range("K113").Value = 3999.1428
$ws.Range("M113").Value = -1829.1428
$ws.Range("H113").Value = 3999.1428
$ws.Range("I113").Value = 3999.1428
# row 134
$ws.Range("I134").Value = 1643.3422
$ws.Range("L134").Value = 33939
$ws.Range("K134").Value = 4930.0266
$ws.Range("J134").Value = 11313
$ws.Range("H134").Value = 4380.0376
$ws.Range("M134").Value = -2395.0266
$ws.Range("N134").Value = -39009

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 62
$ws.Range("K62").Value = 13893969
$ws.Range("M62").Value = -13893345
$ws.Range("H62").Value = 12509572
$ws.Range("I62").Value = 13893969
# row 65
$ws.Range("K65").Value = 69469845
$ws.Range("M65").Value = -69466725
$ws.Range("H65").Value = 12509572
$ws.Range("I65").Value = 13893969
# row 76
$ws.Range("I76").Value = 5319.25
$ws.Range("M76").Value = -5004.25
$ws.Range("H76").Value = 5319.25
$ws.Range("K76").Value = 5319.25
# row 79
$ws.Range("I79").Value = 5319.25
$ws.Range("H79").Value = 5319.25
$ws.Range("K79").Value = 5319.25
$ws.Range("M79").Value = -4227.25
# row 105
$ws.Range("H105").Value = 5953710.5
$ws.Range("J105").Value = 4996
$ws.Range("L105").Value = 4996
$ws.Range("I105").Value = 7143453
$ws.Range("M105").Value = -7141706
$ws.Range("N105").Value = -8490
$ws.Range("K105").Value = 7143453
# row 122
$ws.Range("M122").Value = -2026.299999999999
$ws.Range("K122").Value = 4476.299999999999
$ws.Range("H122").Value = 2191.5625
$ws.Range("I122").Value = 1492.1
# row 132
$ws.Range("J132").Value = 10538.462
$ws.Range("N132").Value = -36675.386
$ws.Range("H132").Value = 7481.2144
$ws.Range("L132").Value = 31615.386
$ws.Range("I132").Value = 4831.6
$ws.Range("M132").Value = -11964.8
$ws.Range("K132").Value = 14494.8

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 16
$ws.Range("L16").Value = 8250
$ws.Range("J16").Value = 2750
$ws.Range("N16").Value = -8596
$ws.Range("I16").Value = 672.5
$ws.Range("M16").Value = -1844.5
$ws.Range("H16").Value = 1365
$ws.Range("K16").Value = 2017.5
# row 34
$ws.Range("N34").Value = -19384.059
$ws.Range("H34").Value = 5456.35
$ws.Range("L34").Value = 19216.059
$ws.Range("J34").Value = 6405.353
# row 75
$ws.Range("J75").Value = 37045384
$ws.Range("L75").Value = 111136152
$ws.Range("H75").Value = 55561876
$ws.Range("N75").Value = -111138148
# row 78
$ws.Range("N78").Value = -333418440
$ws.Range("L78").Value = 333408456
$ws.Range("H78").Value = 55561876
$ws.Range("J78").Value = 37045384
# row 107
$ws.Range("N107").Value = -6589.6875
$ws.Range("J107").Value = 916.5625
$ws.Range("L107").Value = 2749.6875
$ws.Range("H107").Value = 862.6316
# row 126
$ws.Range("H126").Value = 1630
$ws.Range("M126").Value = 50
$ws.Range("K126").Value = 4890
$ws.Range("I126").Value = 1630

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 97
$ws.Range("K97").Value = 1385.1428
$ws.Range("I97").Value = 1385.1428
$ws.Range("H97").Value = 1366.25
$ws.Range("M97").Value = -889.1428000000001
# row 132
$ws.Range("J132").Value = 11605.571
$ws.Range("N132").Value = -39876.713
$ws.Range("H132").Value = 4986.077
$ws.Range("L132").Value = 34816.713
$ws.Range("I132").Value = 2547.3157
$ws.Range("M132").Value = -5111.9471
$ws.Range("K132").Value = 7641.9471

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Range("J22").Value = 6000.5
$ws.Range("M22").Value = -265
$ws.Range("H22").Value = 1648.1
$ws.Range("I22").Value = 560
$ws.Range("L22").Value = 6000.5
$ws.Range("K22").Value = 560
$ws.Range("N22").Value = -6590.5
# row 27
$ws.Range("L27").Value = 6000.5
$ws.Range("M27").Value = -453
$ws.Range("I27").Value = 560
$ws.Range("J27").Value = 6000.5
$ws.Range("H27").Value = 1648.1
$ws.Range("K27").Value = 560
$ws.Range("N27").Value = -6214.5
# row 68
$ws.Range("L68").Value = 5837.615
$ws.Range("J68").Value = 5837.615
$ws.Range("H68").Value = 4492.591
$ws.Range("N68").Value = -7335.615
# row 71
$ws.Range("N71").Value = -36676.075
$ws.Range("H71").Value = 4492.591
$ws.Range("L71").Value = 29188.075
$ws.Range("J71").Value = 5837.615
# row 100
$ws.Range("I100").Value = 3514.7273
$ws.Range("M100").Value = -2973.7273
$ws.Range("K100").Value = 3514.7273
$ws.Range("H100").Value = 4584.3687

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 4
$ws.Range("N4").Value = -489.66666
$ws.Range("K4").Value = 90.5
$ws.Range("H4").Value = 130.46153
$ws.Range("I4").Value = 90.5
$ws.Range("M4").Value = 22.5
$ws.Range("L4").Value = 263.66666
$ws.Range("J4").Value = 263.66666
# row 100
$ws.Range("L100").Value = 2196.25
$ws.Range("J100").Value = 1098.125
$ws.Range("N100").Value = -3278.25
$ws.Range("H100").Value = 737.3333
# row 122
$ws.Range("M122").Value = -3694.777900000001
$ws.Range("K122").Value = 6144.777900000001
$ws.Range("H122").Value = 2649.0513
$ws.Range("I122").Value = 2048.2593
# row 135
$ws.Range("L135").Value = 73000
$ws.Range("N135").Value = -83140
$ws.Range("J135").Value = 73000
$ws.Range("H135").Value = 73000
# row 136
$ws.Range("K136").Value = 4738.857
$ws.Range("J136").Value = 6293.5
$ws.Range("H136").Value = 3991.372
$ws.Range("I136").Value = 1579.619
$ws.Range("N136").Value = -23980.5
$ws.Range("L136").Value = 18880.5
$ws.Range("M136").Value = -2188.857
